$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("R13").Value
Write-Host $v
